# Fixes the two teams' batting/bowling scorecards that were mixed up
# (tournament closure exception handled): corrects the left block
# (Rajastan Australia) and right block (Punjab Pakistan) entries on the
# "Validation" sheet, row by row.
#
# Note: a handful of cells hold over-counts such as "2.0"/"12.0" stored
# as literal text (not numbers). Assigning a bare numeric-looking string
# via COM auto-converts it to a number and drops the trailing zero, so
# those are written with a leading apostrophe to force text semantics,
# matching the original workbook's intent.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'Mohammad Rizwan'
$ws.Range("B2").Value = 10
$ws.Range("C2").Value = 3
$ws.Range("D2").Value = 'Caught'
$ws.Range("E2").Value = ' Josh Hazlewood'
$ws.Range("J2").Value = 'David Warner'
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 1
$ws.Range("N2").Value = ' Haris Rauf'
$ws.Range("A3").Value = 'Babar Azam(C)'
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 'LBW'
$ws.Range("E3").Value = ' Josh Hazlewood'
$ws.Range("J3").Value = 'Aaron Finch(C)'
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 'Caught'
$ws.Range("N3").Value = ' Haris Rauf'
$ws.Range("A4").Value = 'Fakhar Zaman'
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 'LBW'
$ws.Range("E4").Value = ' Josh Hazlewood'
$ws.Range("J4").Value = 'Mitchell Marsh'
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 1
$ws.Range("N4").Value = ' Haris Rauf'
$ws.Range("A5").Value = 'Mohammad Hafeez'
$ws.Range("B5").Value = 13
$ws.Range("C5").Value = 7
$ws.Range("E5").Value = ' Mitchell Starc'
$ws.Range("J5").Value = 'Steve Smith'
$ws.Range("K5").Value = 83
$ws.Range("L5").Value = 25
$ws.Range("M5").Value = 'Bowled'
$ws.Range("N5").Value = ' Shaheen Afridi'
$ws.Range("A6").Value = 'Shoaib Malik'
$ws.Range("B6").Value = 13
$ws.Range("D6").Value = 'Bowled'
$ws.Range("E6").Value = ' Josh Hazlewood'
$ws.Range("J6").Value = 'Glenn Maxwell'
$ws.Range("K6").Value = 16
$ws.Range("L6").Value = 5
$ws.Range("M6").Value = 'LBW'
$ws.Range("N6").Value = ' Hasan Ali'
$ws.Range("A7").Value = 'Asif Ali'
$ws.Range("B7").Value = 4
$ws.Range("C7").Value = 3
$ws.Range("D7").Value = 'Caught'
$ws.Range("E7").Value = ' Pat Cummins'
$ws.Range("J7").Value = 'Matthew Wade'
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 1
$ws.Range("N7").Value = ' Hasan Ali'
$ws.Range("A8").Value = 'Shadab Khan'
$ws.Range("B8").Value = 11
$ws.Range("C8").Value = 3
$ws.Range("D8").Value = 'Bowled'
$ws.Range("E8").Value = ' Pat Cummins'
$ws.Range("J8").Value = 'Marcus Stionis'
$ws.Range("K8").Value = 17
$ws.Range("L8").Value = 9
$ws.Range("M8").Value = 'Bowled'
$ws.Range("N8").Value = ' Hasan Ali'
$ws.Range("A9").Value = 'Imad Wasim'
$ws.Range("B9").Value = 49
$ws.Range("C9").Value = 17
$ws.Range("E9").Value = ' Mitchell Starc'
$ws.Range("J9").Value = 'Pat Cummins'
$ws.Range("K9").Value = 32
$ws.Range("L9").Value = 10
$ws.Range("M9").Value = '* NOT OUT'
$ws.Range("N9").Value = ' '
$ws.Range("A10").Value = 'Hasan Ali'
$ws.Range("B10").Value = 38
$ws.Range("C10").Value = 15
$ws.Range("E10").Value = ' Marcus Stionis'
$ws.Range("J10").Value = 'Mitchell Starc'
$ws.Range("K10").Value = 36
$ws.Range("L10").Value = 15
$ws.Range("M10").Value = 'NOT OUT'
$ws.Range("N10").Value = ' '
$ws.Range("A11").Value = 'Shaheen Afridi'
$ws.Range("B11").Value = 11
$ws.Range("C11").Value = 6
$ws.Range("D11").Value = 'Bowled'
$ws.Range("E11").Value = ' Adam Zampa'
$ws.Range("J11").Value = 'Adam Zampa'
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = ' '
$ws.Range("N11").Value = ' '
$ws.Range("A12").Value = 'Haris Rauf'
$ws.Range("B12").Value = 33
$ws.Range("C12").Value = 11
$ws.Range("D12").Value = 'NOT OUT'
$ws.Range("E12").Value = ' '
$ws.Range("J12").Value = 'Josh Hazlewood'
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = ' '
$ws.Range("N12").Value = ' '
$ws.Range("A16").Value = 182
$ws.Range("C16").Value = "'12.0"
$ws.Range("D16").Value = 72
$ws.Range("J16").Value = 184
$ws.Range("K16").Value = 7
$ws.Range("L16").Value = "'11.2"
$ws.Range("M16").Value = 68
$ws.Range("A21").Value = 'Mitchell Starc'
$ws.Range("B21").Value = "'2.0"
$ws.Range("C21").Value = 31
$ws.Range("D21").Value = 2
$ws.Range("E21").Value = 15.5
$ws.Range("J21").Value = 'Hasan Ali'
$ws.Range("K21").Value = "'2.0"
$ws.Range("L21").Value = 23
$ws.Range("M21").Value = 3
$ws.Range("N21").Value = 11.5
$ws.Range("A22").Value = 'Pat Cummins'
$ws.Range("B22").Value = "'2.0"
$ws.Range("C22").Value = 38
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = 19
$ws.Range("J22").Value = 'Imad Wasim'
$ws.Range("K22").Value = "'2.0"
$ws.Range("L22").Value = 29
$ws.Range("M22").Value = 0
$ws.Range("N22").Value = 14.5
$ws.Range("A23").Value = 'Marcus Stionis'
$ws.Range("B23").Value = "'2.0"
$ws.Range("C23").Value = 30
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 15
$ws.Range("J23").Value = 'Shadab Khan'
$ws.Range("K23").Value = "'2.0"
$ws.Range("L23").Value = 34
$ws.Range("M23").Value = 0
$ws.Range("N23").Value = 17
$ws.Range("A24").Value = 'Josh Hazlewood'
$ws.Range("B24").Value = "'3.0"
$ws.Range("C24").Value = 43
$ws.Range("D24").Value = 4
$ws.Range("E24").Value = 14.33
$ws.Range("J24").Value = 'Haris Rauf'
$ws.Range("K24").Value = "'3.0"
$ws.Range("L24").Value = 57
$ws.Range("M24").Value = 3
$ws.Range("N24").Value = 19
$ws.Range("A25").Value = 'Adam Zampa'
$ws.Range("B25").Value = "'3.0"
$ws.Range("C25").Value = 40
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 13.33
$ws.Range("J25").Value = 'Shaheen Afridi'
$ws.Range("K25").Value = "'2.2"
$ws.Range("L25").Value = 41
$ws.Range("M25").Value = 1
$ws.Range("N25").Value = 18.64
